# Update the "log" workbook:
#  - Switch active/selected sheet from Sheet2 back to Sheet1
#  - On Sheet1: select cell D3, and widen column B
#  - On Sheet2: select cell D25

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Widen column B on Sheet1 (~30.875 character units once stored in the
# workbook's column metadata).
$ws1.Columns.Item(2).ColumnWidth = 30.142857142857142

# Set selection on Sheet2 first (while it's still active) so the saved
# selection for that sheet is D25, then leave Sheet1 as the active sheet.
$ws2.Activate()
$ws2.Range("D25").Select()

$ws1.Activate()
$ws1.Range("D3").Select()
